$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1285113333333333
$ws.Range("H2").Value = 0.385534
$ws.Range("I2").Value = 0.03749201237720504
$ws.Range("J2").Value = 0.03749201237720504
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.110264333333333
$ws.Range("N2").Value = 6.330793
$ws.Range("O2").Value = 0.3832041185227171
$ws.Range("P2").Value = 0.3832041185227171
$ws.Range("Q2").Value = 0.2711928831624444
$ws.Range("R2").Value = 2.440735948462
$ws.Range("S2").Value = 0.01436709355464965
$ws.Range("T2").Value = 0.01436709355464966
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1285113333333333
$ws.Range("H3").Value = 0.385534
$ws.Range("I3").Value = 0.03749201237720504
$ws.Range("J3").Value = 0.03749201237720504
$ws.Range("O3").Value = 0.0946183755984393
$ws.Range("P3").Value = 0.0946183755984393
$ws.Range("Q3").Value = 0.06696125860444443
$ws.Range("R3").Value = 0.6026513274399999
$ws.Range("S3").Value = 0.003547433309047721
$ws.Range("T3").Value = 0.003547433309047722
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1285113333333333
$ws.Range("H4").Value = 0.385534
$ws.Range("I4").Value = 0.03749201237720504
$ws.Range("J4").Value = 0.03749201237720504
$ws.Range("M4").Value = 2.065388333333333
$ws.Range("N4").Value = 6.196165
$ws.Range("O4").Value = 0.3750550597762889
$ws.Range("P4").Value = 0.3750550597762889
$ws.Range("Q4").Value = 0.2654258085677778
$ws.Range("R4").Value = 2.38883227711
$ws.Range("S4").Value = 0.014061568943266
$ws.Range("T4").Value = 0.014061568943266
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1285113333333333
$ws.Range("H5").Value = 0.385534
$ws.Range("I5").Value = 0.03749201237720504
$ws.Range("J5").Value = 0.03749201237720504
$ws.Range("M5").Value = 0.8101876666666666
$ws.Range("N5").Value = 2.430563
$ws.Range("O5").Value = 0.1471224461025547
$ws.Range("P5").Value = 0.1471224461025547
$ws.Range("Q5").Value = 0.1041182972935556
$ws.Range("R5").Value = 0.9370646756419999
$ws.Range("S5").Value = 0.005515916570241663
$ws.Range("T5").Value = 0.005515916570241664
$ws.Range("I6").Value = 0.7552862722193517
$ws.Range("J6").Value = 0.755286272219352
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.110264333333333
$ws.Range("N6").Value = 6.330793
$ws.Range("O6").Value = 0.3832041185227171
$ws.Range("P6").Value = 0.3832041185227171
$ws.Range("Q6").Value = 5.463250670980666
$ws.Range("R6").Value = 49.16925603882599
$ws.Range("S6").Value = 0.2894288101781256
$ws.Range("T6").Value = 0.2894288101781257
$ws.Range("I7").Value = 0.7552862722193517
$ws.Range("J7").Value = 0.755286272219352
$ws.Range("O7").Value = 0.0946183755984393
$ws.Range("P7").Value = 0.0946183755984393
$ws.Range("S7").Value = 0.0714639601891957
$ws.Range("T7").Value = 0.07146396018919572
$ws.Range("I8").Value = 0.7552862722193517
$ws.Range("J8").Value = 0.755286272219352
$ws.Range("M8").Value = 2.065388333333333
$ws.Range("N8").Value = 6.196165
$ws.Range("O8").Value = 0.3750550597762889
$ws.Range("P8").Value = 0.3750550597762889
$ws.Range("Q8").Value = 5.347071463836667
$ws.Range("R8").Value = 48.12364317453
$ws.Range("S8").Value = 0.2832739379754394
$ws.Range("T8").Value = 0.2832739379754394
$ws.Range("I9").Value = 0.7552862722193517
$ws.Range("J9").Value = 0.755286272219352
$ws.Range("M9").Value = 0.8101876666666666
$ws.Range("N9").Value = 2.430563
$ws.Range("O9").Value = 0.1471224461025547
$ws.Range("P9").Value = 0.1471224461025547
$ws.Range("Q9").Value = 2.097489989107333
$ws.Range("R9").Value = 18.877409901966
$ws.Range("S9").Value = 0.1111195638765911
$ws.Range("T9").Value = 0.1111195638765911
$ws.Range("G10").Value = 0.692415
$ws.Range("H10").Value = 2.077245
$ws.Range("I10").Value = 0.2020057770533527
$ws.Range("J10").Value = 0.2020057770533527
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.110264333333333
$ws.Range("N10").Value = 6.330793
$ws.Range("O10").Value = 0.3832041185227171
$ws.Range("P10").Value = 0.3832041185227171
$ws.Range("Q10").Value = 1.461178678365
$ws.Range("R10").Value = 13.150608105285
$ws.Range("S10").Value = 0.07740944573222654
$ws.Range("T10").Value = 0.07740944573222654
$ws.Range("G11").Value = 0.692415
$ws.Range("H11").Value = 2.077245
$ws.Range("I11").Value = 0.2020057770533527
$ws.Range("J11").Value = 0.2020057770533527
$ws.Range("O11").Value = 0.0946183755984393
$ws.Range("P11").Value = 0.0946183755984393
$ws.Range("Q11").Value = 0.3607851437999999
$ws.Range("R11").Value = 3.2470662942
$ws.Range("S11").Value = 0.01911345848628872
$ws.Range("T11").Value = 0.01911345848628872
$ws.Range("G12").Value = 0.692415
$ws.Range("H12").Value = 2.077245
$ws.Range("I12").Value = 0.2020057770533527
$ws.Range("J12").Value = 0.2020057770533527
$ws.Range("M12").Value = 2.065388333333333
$ws.Range("N12").Value = 6.196165
$ws.Range("O12").Value = 0.3750550597762889
$ws.Range("P12").Value = 0.3750550597762889
$ws.Range("Q12").Value = 1.430105862825
$ws.Range("R12").Value = 12.870952765425
$ws.Range("S12").Value = 0.07576328878790088
$ws.Range("T12").Value = 0.0757632887879009
$ws.Range("G13").Value = 0.692415
$ws.Range("H13").Value = 2.077245
$ws.Range("I13").Value = 0.2020057770533527
$ws.Range("J13").Value = 0.2020057770533527
$ws.Range("M13").Value = 0.8101876666666666
$ws.Range("N13").Value = 2.430563
$ws.Range("O13").Value = 0.1471224461025547
$ws.Range("P13").Value = 0.1471224461025547
$ws.Range("Q13").Value = 0.560986093215
$ws.Range("R13").Value = 5.048874838934999
$ws.Range("S13").Value = 0.02971958404693657
$ws.Range("T13").Value = 0.02971958404693658
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.01787866666666667
$ws.Range("H14").Value = 0.053636
$ws.Range("I14").Value = 0.005215938350090445
$ws.Range("J14").Value = 0.005215938350090446
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.110264333333333
$ws.Range("N14").Value = 6.330793
$ws.Range("O14").Value = 0.3832041185227171
$ws.Range("P14").Value = 0.3832041185227171
$ws.Range("Q14").Value = 0.03772871259422222
$ws.Range("R14").Value = 0.339558413348
$ws.Range("S14").Value = 0.001998769057715244
$ws.Range("T14").Value = 0.001998769057715245
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.01787866666666667
$ws.Range("H15").Value = 0.053636
$ws.Range("I15").Value = 0.005215938350090445
$ws.Range("J15").Value = 0.005215938350090446
$ws.Range("O15").Value = 0.0946183755984393
$ws.Range("P15").Value = 0.0946183755984393
$ws.Range("Q15").Value = 0.009315738862222221
$ws.Range("R15").Value = 0.08384164976
$ws.Range("S15").Value = 0.0004935236139071615
$ws.Range("T15").Value = 0.0004935236139071616
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.01787866666666667
$ws.Range("H16").Value = 0.053636
$ws.Range("I16").Value = 0.005215938350090445
$ws.Range("J16").Value = 0.005215938350090446
$ws.Range("M16").Value = 2.065388333333333
$ws.Range("N16").Value = 6.196165
$ws.Range("O16").Value = 0.3750550597762889
$ws.Range("P16").Value = 0.3750550597762889
$ws.Range("Q16").Value = 0.03692638954888889
$ws.Range("R16").Value = 0.33233750594
$ws.Range("S16").Value = 0.00195626406968261
$ws.Range("T16").Value = 0.00195626406968261
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.01787866666666667
$ws.Range("H17").Value = 0.053636
$ws.Range("I17").Value = 0.005215938350090445
$ws.Range("J17").Value = 0.005215938350090446
$ws.Range("M17").Value = 0.8101876666666666
$ws.Range("N17").Value = 2.430563
$ws.Range("O17").Value = 0.1471224461025547
$ws.Range("P17").Value = 0.1471224461025547
$ws.Range("Q17").Value = 0.01448507522977778
$ws.Range("R17").Value = 0.130365677068
$ws.Range("S17").Value = 0.0007673816087854298
$ws.Range("T17").Value = 0.0007673816087854299
Write-Host "Applied TPM data update to $($wb.ActiveSheet.Name)"
